$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 11 data rows (rows 2-12), shifting remaining data up
$ws.Range("A2:C12").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Append 10 new data rows (rows 12-21) with newly generated samples
$ws.Range("A12").Value = -3.011718273162842
$ws.Range("B12").Value = 2.93963623046875
$ws.Range("C12").Value = 0.5337435007095337
$ws.Range("A13").Value = -0.4014911949634552
$ws.Range("B13").Value = 0.9065240025520324
$ws.Range("C13").Value = -0.6395758986473083
$ws.Range("A14").Value = -0.99250328540802
$ws.Range("B14").Value = 2.58121109008789
$ws.Range("C14").Value = -0.3058907687664032
$ws.Range("A15").Value = -2.1601722240448
$ws.Range("B15").Value = 2.302809238433838
$ws.Range("C15").Value = -0.3118467032909393
$ws.Range("A16").Value = -0.5739079117774963
$ws.Range("B16").Value = -1.005178809165955
$ws.Range("C16").Value = 0.2600758671760559
$ws.Range("A17").Value = -1.500132322311401
$ws.Range("B17").Value = 1.276097536087036
$ws.Range("C17").Value = -1.762957096099854
$ws.Range("A18").Value = -1.534646153450012
$ws.Range("B18").Value = -1.36818540096283
$ws.Range("C18").Value = 0.8046622276306152
$ws.Range("A19").Value = -1.083674907684326
$ws.Range("B19").Value = 1.494787216186523
$ws.Range("C19").Value = 0.7244861721992493
$ws.Range("A20").Value = 1.442405581474304
$ws.Range("B20").Value = -1.076344609260559
$ws.Range("C20").Value = 0.901637077331543
$ws.Range("A21").Value = 4.085466861724854
$ws.Range("B21").Value = -3.372739791870117
$ws.Range("C21").Value = 4.743368625640869
